# Weekly update: a new daily price record was inserted for
# "Mapocho Venta Directa de Santiago - Pepino dulce" at row 89, pushing
# all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 89, shifting rows 89:137 down to 90:138.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new record's data.
$ws.Range('A89').Value = 12
$ws.Range('B89').Value = 'Mapocho Venta Directa de Santiago'
$ws.Range('C89').Value = 'Metropolitana'
$ws.Range('D89').Value = 44466
$ws.Range('E89').Value = 13
$ws.Range('F89').Value = 100112043
$ws.Range('G89').Value = 'Pepino dulce'
$ws.Range('H89').Value = 'Cultivar IV Región'
$ws.Range('I89').Value = 'Primera'
$ws.Range('J89').Value = 180
$ws.Range('K89').Value = 20000
$ws.Range('L89').Value = 20000
$ws.Range('M89').Value = 20000
$ws.Range('N89').Value = '$/caja 18 kilos'
$ws.Range('O89').Value = 'Provincia de Limarí'
$ws.Range('P89').Value = 1111
$ws.Range('Q89').Value = 18
$ws.Range('R89').Value = 'Hortaliza'
